# Apply updated symbol list values (Price / Volume(1h)) per commit.
# Source values are stored as literal text in the sheet (t="inlineStr"),
# so we prefix assignments with a leading apostrophe to force Excel to
# store them as text rather than reinterpreting as numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'301.41"
$ws.Range("E2").Value = "'-6.26%"
$ws.Range("D3").Value = "'34.98"
$ws.Range("E3").Value = "'-3.48%"
$ws.Range("D4").Value = "'5.025"
$ws.Range("E4").Value = "'-2.00%"
$ws.Range("D5").Value = "'0.07905"
$ws.Range("E5").Value = "'-2.44%"
$ws.Range("D6").Value = "'1.947"
$ws.Range("E6").Value = "'-9.99%"
$ws.Range("D7").Value = "'7.750"
$ws.Range("E7").Value = "'-3.22%"
$ws.Range("D8").Value = "'4.023"
$ws.Range("E8").Value = "'-2.73%"
$ws.Range("E9").Value = "'5.85%"
$ws.Range("D10").Value = "'0.9227"
$ws.Range("E10").Value = "'-0.45%"
$ws.Range("D11").Value = "'0.1183"
$ws.Range("E11").Value = "'18.54%"
$ws.Range("D12").Value = "'0.1829"
$ws.Range("E12").Value = "'-3.16%"
$ws.Range("D13").Value = "'0.09266"
$ws.Range("E13").Value = "'0.71%"
$ws.Range("E14").Value = "'-1.71%"
$ws.Range("D15").Value = "'0.09871"
$ws.Range("E15").Value = "'-0.59%"
$ws.Range("D16").Value = "'0.001394"
$ws.Range("E16").Value = "'-2.68%"
$ws.Range("D17").Value = "'0.005837"
$ws.Range("E17").Value = "'3.08%"
$ws.Range("D18").Value = "'3.504"
$ws.Range("E18").Value = "'1.42%"
$ws.Range("D19").Value = "'0.3444"
$ws.Range("E19").Value = "'2.12%"
$ws.Range("D20").Value = "'0.1309"
$ws.Range("E20").Value = "'-1.51%"
$ws.Range("D21").Value = "'5.032"
$ws.Range("E21").Value = "'-0.53%"
$ws.Range("E23").Value = "'-2.29%"
$ws.Range("D24").Value = "'0.001213"
$ws.Range("E24").Value = "'-2.41%"
$ws.Range("E25").Value = "'-3.62%"
$ws.Range("E26").Value = "'-3.93%"
$ws.Range("E27").Value = "'-6.78%"
$ws.Range("D39").Value = "'0.01893"
$ws.Range("E39").Value = "'-7.37%"
$ws.Range("D40").Value = "'0.04693"
$ws.Range("E40").Value = "'-5.91%"
$ws.Range("D41").Value = "'0.007573"
$ws.Range("E41").Value = "'-3.11%"
$ws.Range("D42").Value = "'0.009560"
$ws.Range("E42").Value = "'22.22%"
$ws.Range("D43").Value = "'0.1323"
$ws.Range("E43").Value = "'-5.38%"
$ws.Range("D44").Value = "'0.002109"
$ws.Range("E44").Value = "'1.34%"
$ws.Range("D45").Value = "'0.01115"
$ws.Range("E45").Value = "'-8.08%"
$ws.Range("D46").Value = "'0.00005998"
$ws.Range("E46").Value = "'-6.54%"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("E49").Value = "'-31.32%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("E51").Value = "'0.07%"
